$d = $word.ActiveDocument

# Helper: set a paragraph's visible text without touching its trailing
# paragraph mark (doing so on a Paragraph.Range directly can duplicate the
# old tail text in this host, so the end is pulled back one character first).
function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

# Helper: find the (1-based) paragraph index whose text equals $text
# (paragraph text always ends with a trailing CR, hence the `r`).
function Find-ParaIndex([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -eq ($text + "`r")) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Model Description paragraph: mention that routing uses the CPLEX
#    solver right after the distance-minimisation sentence, and replace the
#    old "parcels delivered ... CPLEX solver" sentences with the new
#    queueing-simulation description.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    " At each delivery point, the robot waits for a ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " The routing is accomplished with the aid of the CPLEX solver. At each delivery point, the robot waits for a ",
    2) | Out-Null

$d.Content.Find.Execute(
    "The total number of parcels delivered at each node is a function of the number of parcels to be delivered, the waiting time at each delivery point, and the time it takes to deliver a single parcel. The model uses the CPLEX solver.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "At each delivery point while the robot is waiting, we simulate the operation of a single-server queue in order to determine the percentage of customers served, the average waiting time for served customers, and the average queue length. ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Input-parameters bullet list (numId=3): add a new bullet for the
#    customer arrival rate right before "The waiting time ... delivery
#    point", and extend that bullet with "(also known as the service
#    time)".
# ---------------------------------------------------------------------------

$waitIndex = Find-ParaIndex "The waiting time at each delivery point"
$d.Paragraphs.Item($waitIndex).Range.InsertParagraphBefore()
Set-ParaText $d.Paragraphs.Item($waitIndex) "The customer arrival rate (per minute) at each delivery point"
Set-ParaText $d.Paragraphs.Item($waitIndex + 1) "The waiting time (also known as the service time) at each delivery point"

# ---------------------------------------------------------------------------
# 3) Outputs bullet list (numId=4): reword two bullets to scope them "at
#    each delivery point", and add two new bullets ("average waiting time"
#    / "average queue length") right before the "total emissions" bullet.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "The number of parcels delivered (based on the waiting time and the average delivery time)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The number of parcels delivered at each delivery point",
    2) | Out-Null

$d.Content.Find.Execute(
    "The percentage of parcels successfully delivered",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The percentage of parcels successfully delivered at each delivery point",
    2) | Out-Null

$emissionsIndex = Find-ParaIndex "The total emissions (in gCO2eq) during the delivery operation"
$d.Paragraphs.Item($emissionsIndex).Range.InsertParagraphBefore()
Set-ParaText $d.Paragraphs.Item($emissionsIndex) "The average waiting time at each delivery point"
$emissionsIndex = Find-ParaIndex "The total emissions (in gCO2eq) during the delivery operation"
$d.Paragraphs.Item($emissionsIndex).Range.InsertParagraphBefore()
Set-ParaText $d.Paragraphs.Item($emissionsIndex) "The average queue length at each delivery point"

Write-Host "Done"
